$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(31, 8).Value = 1767.7333
$ws.Cells.Item(31, 9).Value = 401.6
$ws.Cells.Item(31, 10).Value = 4500
$ws.Cells.Item(31, 11).Value = 1204.8
$ws.Cells.Item(31, 12).Value = 13500
$ws.Cells.Item(31, 13).Value = -974.8000000000002
$ws.Cells.Item(31, 14).Value = -13960
$ws.Cells.Item(32, 8).Value = 71429144
$ws.Cells.Item(32, 9).Value = 166667040
$ws.Cells.Item(32, 10).Value = 730.375
$ws.Cells.Item(32, 11).Value = 166667040
$ws.Cells.Item(32, 12).Value = 730.375
$ws.Cells.Item(32, 13).Value = -166666714
$ws.Cells.Item(32, 14).Value = -1382.375
$ws.Cells.Item(38, 8).Value = 2723.5
$ws.Cells.Item(38, 9).Value = 157.85715
$ws.Cells.Item(38, 10).Value = 4356.1816
$ws.Cells.Item(38, 11).Value = 473.57145
$ws.Cells.Item(38, 12).Value = 13068.5448
$ws.Cells.Item(38, 13).Value = -101.57145
$ws.Cells.Item(38, 14).Value = -13812.5448
$ws.Cells.Item(39, 8).Value = 8772181
$ws.Cells.Item(39, 10).Value = 18518828
$ws.Cells.Item(39, 12).Value = 55556484
$ws.Cells.Item(39, 14).Value = -55557076
$ws.Cells.Item(70, 8).Value = 1841
$ws.Cells.Item(70, 9).Value = 1498.5
$ws.Cells.Item(70, 10).Value = 1917.1111
$ws.Cells.Item(70, 11).Value = 4495.5
$ws.Cells.Item(70, 12).Value = 5751.3333
$ws.Cells.Item(70, 13).Value = -4225.5
$ws.Cells.Item(70, 14).Value = -6291.3333
$ws.Cells.Item(73, 8).Value = 1841
$ws.Cells.Item(73, 9).Value = 1498.5
$ws.Cells.Item(73, 10).Value = 1917.1111
$ws.Cells.Item(73, 11).Value = 4495.5
$ws.Cells.Item(73, 12).Value = 5751.3333
$ws.Cells.Item(73, 13).Value = -3559.5
$ws.Cells.Item(73, 14).Value = -7623.3333
$ws.Cells.Item(119, 8).Value = 5966.6665
$ws.Cells.Item(119, 10).Value = 5966.6665
$ws.Cells.Item(119, 12).Value = 17899.9995
$ws.Cells.Item(119, 14).Value = -27575.9995
$ws.Cells.Item(129, 8).Value = 1314.2858
$ws.Cells.Item(129, 9).Value = 425
$ws.Cells.Item(129, 10).Value = 2500
$ws.Cells.Item(129, 11).Value = 1275
$ws.Cells.Item(129, 12).Value = 7500
$ws.Cells.Item(129, 13).Value = 3725
$ws.Cells.Item(129, 14).Value = -17500
$ws.Cells.Item(135, 8).Value = 39751.81
$ws.Cells.Item(135, 9).Value = 49031.477
$ws.Cells.Item(135, 10).Value = 777.2
$ws.Cells.Item(135, 11).Value = 441283.293
$ws.Cells.Item(135, 12).Value = 6994.8
$ws.Cells.Item(135, 13).Value = -438748.293
$ws.Cells.Item(135, 14).Value = -12064.8
$ws.Cells.Item(138, 8).Value = 2348.662
$ws.Cells.Item(138, 9).Value = 1993
$ws.Cells.Item(138, 10).Value = 2462.9822
$ws.Cells.Item(138, 11).Value = 5979
$ws.Cells.Item(138, 12).Value = 7388.946599999999
$ws.Cells.Item(138, 13).Value = -839
$ws.Cells.Item(138, 14).Value = -17668.9466

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 2075.0425
$ws.Cells.Item(97, 9).Value = 1083.6774
$ws.Cells.Item(97, 10).Value = 3995.8125
$ws.Cells.Item(97, 11).Value = 1083.6774
$ws.Cells.Item(97, 12).Value = 3995.8125
$ws.Cells.Item(97, 13).Value = -587.6774
$ws.Cells.Item(97, 14).Value = -4987.8125
$ws.Cells.Item(110, 8).Value = 3913.1765
$ws.Cells.Item(110, 9).Value = 5044.5
$ws.Cells.Item(110, 10).Value = 1198
$ws.Cells.Item(110, 11).Value = 5044.5
$ws.Cells.Item(110, 12).Value = 1198
$ws.Cells.Item(110, 13).Value = -2999.5
$ws.Cells.Item(110, 14).Value = -5288
$ws.Cells.Item(122, 8).Value = 1582.4286
$ws.Cells.Item(122, 9).Value = 1445.7
$ws.Cells.Item(122, 11).Value = 4337.1
$ws.Cells.Item(122, 13).Value = -1887.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 750.7879
$ws.Cells.Item(94, 9).Value = 664.6818
$ws.Cells.Item(94, 10).Value = 923
$ws.Cells.Item(94, 11).Value = 664.6818
$ws.Cells.Item(94, 12).Value = 923
$ws.Cells.Item(94, 13).Value = -213.6818
$ws.Cells.Item(94, 14).Value = -1825
$ws.Cells.Item(99, 8).Value = 1801.1364
$ws.Cells.Item(99, 9).Value = 1080.2
$ws.Cells.Item(99, 10).Value = 3346
$ws.Cells.Item(99, 11).Value = 1080.2
$ws.Cells.Item(99, 12).Value = 3346
$ws.Cells.Item(99, 13).Value = 417.8
$ws.Cells.Item(99, 14).Value = -6342
$ws.Cells.Item(134, 8).Value = 2775.7917
$ws.Cells.Item(134, 9).Value = 2564.1333
$ws.Cells.Item(134, 10).Value = 3128.5557
$ws.Cells.Item(134, 11).Value = 7692.3999
$ws.Cells.Item(134, 12).Value = 9385.667099999999
$ws.Cells.Item(134, 13).Value = -5157.3999
$ws.Cells.Item(134, 14).Value = -14455.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(5, 8).Value = 1000263
$ws.Cells.Item(5, 10).Value = 1000263
$ws.Cells.Item(5, 12).Value = 1000263
$ws.Cells.Item(5, 14).Value = -1000487
$ws.Cells.Item(31, 8).Value = 3105.3057
$ws.Cells.Item(31, 9).Value = 1678.4642
$ws.Cells.Item(31, 11).Value = 1678.4642
$ws.Cells.Item(31, 13).Value = -1383.4642
$ws.Cells.Item(34, 8).Value = 3105.3057
$ws.Cells.Item(34, 9).Value = 1678.4642
$ws.Cells.Item(34, 11).Value = 1678.4642
$ws.Cells.Item(34, 13).Value = -1476.4642
$ws.Cells.Item(132, 8).Value = 3186.476
$ws.Cells.Item(132, 9).Value = 2757.375
$ws.Cells.Item(132, 10).Value = 4559.6
$ws.Cells.Item(132, 11).Value = 8272.125
$ws.Cells.Item(132, 12).Value = 13678.8
$ws.Cells.Item(132, 13).Value = -5742.125
$ws.Cells.Item(132, 14).Value = -18738.8
$ws.Cells.Item(134, 8).Value = 5280.276
$ws.Cells.Item(134, 9).Value = 6930.8423
$ws.Cells.Item(134, 10).Value = 2144.2
$ws.Cells.Item(134, 11).Value = 20792.5269
$ws.Cells.Item(134, 12).Value = 6432.599999999999
$ws.Cells.Item(134, 13).Value = -18257.5269
$ws.Cells.Item(134, 14).Value = -11502.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(124, 8).Value = 3721.111
$ws.Cells.Item(124, 9).Value = 698
$ws.Cells.Item(124, 10).Value = 7500
$ws.Cells.Item(124, 11).Value = 2094
$ws.Cells.Item(124, 12).Value = 22500
$ws.Cells.Item(124, 13).Value = 2816
$ws.Cells.Item(124, 14).Value = -32320
$ws.Cells.Item(131, 8).Value = 976.6842
$ws.Cells.Item(131, 10).Value = 1052.22
$ws.Cells.Item(131, 12).Value = 3156.66
$ws.Cells.Item(131, 14).Value = -13236.66

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 690.8077
$ws.Cells.Item(107, 9).Value = 488.66666
$ws.Cells.Item(107, 11).Value = 488.66666
$ws.Cells.Item(107, 13).Value = 1431.33334
$ws.Cells.Item(113, 8).Value = 1607.0588
$ws.Cells.Item(113, 9).Value = 1620.4667
$ws.Cells.Item(113, 11).Value = 1620.4667
$ws.Cells.Item(113, 13).Value = 549.5333000000001
$ws.Cells.Item(122, 8).Value = 2944.4546
$ws.Cells.Item(122, 9).Value = 2911.125
$ws.Cells.Item(122, 10).Value = 3033.3333
$ws.Cells.Item(122, 11).Value = 8733.375
$ws.Cells.Item(122, 12).Value = 9099.999899999999
$ws.Cells.Item(122, 13).Value = -6283.375
$ws.Cells.Item(122, 14).Value = -13999.9999
$ws.Cells.Item(132, 8).Value = 3461.25
$ws.Cells.Item(132, 9).Value = 2469.75
$ws.Cells.Item(132, 10).Value = 4948.5
$ws.Cells.Item(132, 11).Value = 7409.25
$ws.Cells.Item(132, 12).Value = 14845.5
$ws.Cells.Item(132, 13).Value = -4879.25
$ws.Cells.Item(132, 14).Value = -19905.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1335.4117
$ws.Cells.Item(46, 9).Value = 1200.1666
$ws.Cells.Item(46, 11).Value = 1200.1666
$ws.Cells.Item(46, 13).Value = -1012.1666
$ws.Cells.Item(132, 8).Value = 3146
$ws.Cells.Item(132, 9).Value = 2433.3333
$ws.Cells.Item(132, 10).Value = 4749.5
$ws.Cells.Item(132, 11).Value = 7299.999899999999
$ws.Cells.Item(132, 12).Value = 14248.5
$ws.Cells.Item(132, 13).Value = -4769.999899999999
$ws.Cells.Item(132, 14).Value = -19308.5
$ws.Cells.Item(133, 8).Value = 48102.832
$ws.Cells.Item(133, 10).Value = 48102.832
$ws.Cells.Item(133, 12).Value = 48102.832
$ws.Cells.Item(133, 14).Value = -53162.832
$ws.Cells.Item(136, 8).Value = 1652.5883
$ws.Cells.Item(136, 9).Value = 1443.375
$ws.Cells.Item(136, 10).Value = 5000
$ws.Cells.Item(136, 11).Value = 4330.125
$ws.Cells.Item(136, 12).Value = 15000
$ws.Cells.Item(136, 13).Value = -1780.125
$ws.Cells.Item(136, 14).Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 5250
$ws.Cells.Item(4, 9).Value = 2000
$ws.Cells.Item(4, 10).Value = 6333.3335
$ws.Cells.Item(4, 11).Value = 2000
$ws.Cells.Item(4, 12).Value = 6333.3335
$ws.Cells.Item(4, 13).Value = -1887
$ws.Cells.Item(4, 14).Value = -6559.3335
$ws.Cells.Item(17, 8).Value = 28876
$ws.Cells.Item(17, 9).Value = 30336
$ws.Cells.Item(17, 10).Value = 28000
$ws.Cells.Item(17, 11).Value = 30336
$ws.Cells.Item(17, 12).Value = 28000
$ws.Cells.Item(17, 13).Value = -30164
$ws.Cells.Item(17, 14).Value = -28344
$ws.Cells.Item(136, 8).Value = 1105.317
$ws.Cells.Item(136, 9).Value = 854.8205
$ws.Cells.Item(136, 11).Value = 2564.4615
$ws.Cells.Item(136, 13).Value = -14.46150000000034
